$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink that was on F2 (mailto:rjaini@vmware.com1)
$ws.Hyperlinks.Delete()

# --- Row 1 headers (text unchanged, but now reference different shared-string slots
#     because the old "...1" demo strings are removed) ---
$ws.Range("A1").Value = "UserId"
$ws.Range("H1").Value = "CreatedDate"

# --- Row 2 data: replace the old demo values with the new ones ---
$ws.Range("B2").Value = 12
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = "rjaini@vmware.com"
$ws.Range("G2").Value = "Admin"
$ws.Range("H2").Value = 44146

# --- Column A: unhide and resize to fit its new "UserId" values ---
$ws.Columns.Item(1).Hidden = $false
$ws.Columns.Item(1).ColumnWidth = 5.666666666666666

# --- Column H: resize for the new (shorter) date values ---
$ws.Columns.Item(8).ColumnWidth = 9.666666666666668

# --- View state: scroll back to A1 and move the selection to B3 ---
$ws.Range("B3").Select()
